$d = $word.ActiveDocument

# Locate the table / cells that hold the Libsodium timing numbers we need
# to update: "29.56" -> "28.69" and "~34K" -> "~35K".
$targetTable = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $tbl = $d.Tables.Item($i)
    if ($tbl.Range.Text -like "*29.56*") {
        $targetTable = $tbl
    }
}

$cell1 = $targetTable.Cell(2, 1)   # currently "29.56"
$cell2 = $targetTable.Cell(2, 2)   # currently "~34K"

# ---------------------------------------------------------------------
# Cell (2,1): "29.56" -> "28.69", ending up split across four runs:
#   "2" | "8" | "." | "69"
# ---------------------------------------------------------------------
$r1 = $cell1.Range
$start1 = $r1.Start
$end1 = $r1.End

$full1 = $d.Range($start1, $end1 - 1)
$full1.Text = "28.69"

# Force run boundaries by toggling a character property on/off - Word
# keeps the run split even once the explicit formatting is removed again.
$split1a = $d.Range($start1, $start1 + 1)          # "2"
$split1a.Font.Bold = 1
$split1a.Font.Bold = 0

$split1b = $d.Range($start1 + 1, $start1 + 2)      # "8"
$split1b.Font.Bold = 1
$split1b.Font.Bold = 0

$split1c = $d.Range($start1 + 2, $start1 + 3)      # "."
$split1c.Font.Bold = 1
$split1c.Font.Bold = 0
# remaining "69" stays as the trailing run

# ---------------------------------------------------------------------
# Cell (2,2): "~34K" -> "~35K", ending up split across three runs:
#   "~3" | "5" | "K"
# ---------------------------------------------------------------------
$r2 = $cell2.Range
$start2 = $r2.Start
$end2 = $r2.End

$full2 = $d.Range($start2, $end2 - 1)
$full2.Text = "~35K"

$split2a = $d.Range($start2, $start2 + 2)          # "~3"
$split2a.Font.Bold = 1
$split2a.Font.Bold = 0

$split2b = $d.Range($start2 + 2, $start2 + 3)      # "5"
$split2b.Font.Bold = 1
$split2b.Font.Bold = 0
# remaining "K" stays as the trailing run

Write-Output "Updated timing cells: 29.56 -> 28.69, ~34K -> ~35K"
